$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 0.34837547187382256
$ws.Range("A2").Value = -0.038899119782232816
$ws.Range("A3").Value = -0.0089999990383127226
$ws.Range("A4").Value = -0.011999999757893676
$ws.Range("A5").Value = -0.0059999990493277977
$ws.Range("A6").Value = -0.02187340111522218
$ws.Range("A7").Value = -0.019999998854618894
$ws.Range("A8").Value = -0.019999998854034473
$ws.Range("A9").Value = -0.0059999990218928545
$ws.Range("A10").Value = -0.0059999990262298297
$ws.Range("A11").Value = 0.046537649854304419
$ws.Range("A12").Value = -0.0059999990255388269
$ws.Range("A13").Value = -0.005999999020949609
$ws.Range("A14").Value = 0.042427068376224675
$ws.Range("A15").Value = -0.0059999990132855174
$ws.Range("A16").Value = -0.031159799092430163
$ws.Range("A17").Value = -0.0059999990015970894
$ws.Range("A18").Value = -0.0089999989653213319
$ws.Range("A19").Value = -0.0089999990668325758
$ws.Range("A20").Value = -0.0089999990585187817
$ws.Range("A21").Value = -0.00899999905729576
$ws.Range("A22").Value = -0.0089999990564555432
$ws.Range("A23").Value = -0.0089999990296432131
$ws.Range("A24").Value = -0.041999998627201052
$ws.Range("A25").Value = -0.041999998619611567
$ws.Range("A26").Value = -0.0059999990169039563
$ws.Range("A27").Value = -0.005999999011366608
$ws.Range("A28").Value = -0.0059999989899175432
$ws.Range("A29").Value = -0.011999998904720144
$ws.Range("A30").Value = -0.019999998803184482
$ws.Range("A31").Value = -0.014999998852150753
$ws.Range("A32").Value = -0.020999998780148132
$ws.Range("A33").Value = -0.0059999989566916767
